$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibition) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 648
$ws1.Range("G2").Value = 70
$ws1.Range("F3").Value = 735
$ws1.Range("F4").Value = 952
$ws1.Range("F5").Value = 734
$ws1.Range("F6").Value = 844
$ws1.Range("F7").Value = 409
$ws1.Range("F8").Value = 618
$ws1.Range("F10").Value = 1224
$ws1.Range("F11").Value = 648
$ws1.Range("F13").Value = 516
$ws1.Range("F15").Value = 16
$ws1.Range("F16").Value = 591
$ws1.Range("F17").Value = 2
$ws1.Range("F18").Value = 362
$ws1.Range("F22").Value = 91
$ws1.Range("F23").Value = 587
$ws1.Range("F25").Value = 807

# ---- Sheet "演出" (Performance) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 328
$ws2.Range("F9").Value = 227
$ws2.Range("F10").Value = 51
$ws2.Range("F13").Value = 101

# ---- Sheet "全部类型" (All Types, combined) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 648
$ws4.Range("G4").Value = 70
$ws4.Range("F6").Value = 328
$ws4.Range("F7").Value = 735
$ws4.Range("F8").Value = 952
$ws4.Range("F9").Value = 734
$ws4.Range("F10").Value = 844
$ws4.Range("F11").Value = 409
$ws4.Range("F12").Value = 618
$ws4.Range("F14").Value = 1224
$ws4.Range("F15").Value = 648
$ws4.Range("F19").Value = 516
$ws4.Range("F22").Value = 16
$ws4.Range("F23").Value = 591
$ws4.Range("F25").Value = 2
$ws4.Range("F26").Value = 362
$ws4.Range("F29").Value = 227
$ws4.Range("F30").Value = 51
$ws4.Range("F34").Value = 101
$ws4.Range("F35").Value = 101
$ws4.Range("F36").Value = 91
$ws4.Range("F37").Value = 587
$ws4.Range("F39").Value = 807
